# Guided Capstone slide - text edits
#
# Shape "Google Shape;36;p1" (the "There are limitations..." textbox):
#   - run 1 loses its trailing "T"
#   - run 2's text is replaced with a longer sentence that leads into
#     "...so the better change is to cut the costs."
#
# Shape "Google Shape;37;p1" (the "There's a suspicion..." textbox):
#   - the single run is split into two runs with new wording and bold styling.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Find-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Name -eq $name) {
            return $sh
        }
    }
    return $null
}

# --- Shape 36: "There are limitations to raise the tickets price. T" / "he better change is to cut the costs." ---
$shape36 = Find-ShapeByName $s "Google Shape;36;p1"
if ($shape36 -eq $null) { $shape36 = $s.Shapes.Item(17) }
$tr36 = $shape36.TextFrame.TextRange

$run1 = $tr36.Runs(1, 1)
$run1.Text = "There are limitations to raise the tickets price. "

$run2 = $tr36.Runs(2, 1)
$run2.Text = "There's a suspicion that Big Mountain is not capitalizing on its facilities as much as it could, so the better change is to cut the costs."

# --- Shape 37: "There's a suspicion that Big Mountain is not capitalizing..." ---
$shape37 = Find-ShapeByName $s "Google Shape;37;p1"
if ($shape37 -eq $null) { $shape37 = $s.Shapes.Item(18) }
$tr37 = $shape37.TextFrame.TextRange

$tr37.Text = "There were not data sources about  how Big Mountain Resort capitalizing on its facilities that we can use to analyzing.  It is just a suspicion."

$part1 = $tr37.Characters(1, 121)
$part1.Font.Size = 10.7
$part1.Font.Bold = -1

$part2 = $tr37.Characters(122, 23)
$part2.Font.Size = 11
$part2.Font.Bold = -1
